$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 0.5087393606160395
$ws.Range("C3").Value = -1.118515468742087
$ws.Range("E3").Value = -0.6296678961043134
$ws.Range("C4").Value = -0.1156872058426073
$ws.Range("E4").Value = -0.5120992642018263
$ws.Range("C5").Value = -0.4084169314491404
$ws.Range("E5").Value = -0.6403426624573716
$ws.Range("C6").Value = -0.2188016966516937
$ws.Range("E6").Value = -0.1561757764150462
$ws.Range("C7").Value = 0.05915234751026066
$ws.Range("E7").Value = 0.04624521867206965
$ws.Range("C8").Value = -0.2979029954603529
$ws.Range("E8").Value = -0.1124510725819206
$ws.Range("C9").Value = 0.07317408757452348
$ws.Range("E9").Value = -0.002181547367274828
$ws.Range("C10").Value = -0.06188089372189953
$ws.Range("E10").Value = -0.07932008107318644
$ws.Range("C11").Value = -0.2199961235931358
$ws.Range("E11").Value = -0.1022879117640763
$ws.Range("C12").Value = 0.0882025545300813
$ws.Range("E12").Value = -0.05573300569792217
$ws.Range("C13").Value = -0.7844010209450802
$ws.Range("E13").Value = -0.3786583343736716
$ws.Range("C14").Value = -0.1808804304865297
$ws.Range("E14").Value = -0.1077309791980285
$ws.Range("C15").Value = -0.8769761459347714
$ws.Range("E15").Value = -0.5354669478056073
$ws.Range("C16").Value = -1.719168896439693
$ws.Range("E16").Value = -0.5989817782328322
$ws.Range("C17").Value = -0.6902657121583777
$ws.Range("E17").Value = -1.134712300966823
$ws.Range("C18").Value = -0.07874066250703748
$ws.Range("E18").Value = -0.1533081888441812
$ws.Range("C19").Value = 0.4305325812036687
$ws.Range("E19").Value = 0.2136583044595852
